# Updates cryptos list values (Price / Volume(1h) columns) per source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the literal as text even
    # when it looks numeric (e.g. "116.49"); Style reset afterwards avoids
    # leaving a stray quote-prefix style applied to the cell.
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "51.648.98"
Set-TextValue "E2" "  +6.28%  "
Set-TextValue "D3" "2.743.64"
Set-TextValue "E3" "  +5.15%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "116.49"
Set-TextValue "D6" "333.29"
Set-TextValue "E6" "  +3.77%  "
Set-TextValue "E7" "  +2.55%  "
Set-TextValue "E8" "  -0.01%  "
Set-TextValue "D9" "0.570"
Set-TextValue "E9" "  +5.86%  "
Set-TextValue "D10" "41.39"
Set-TextValue "E10" "  +5.59%  "
Set-TextValue "D11" "0.0856"
Set-TextValue "E11" "  +6.07%  "
Set-TextValue "D12" "20.12"
Set-TextValue "E12" "  +2.09%  "
Set-TextValue "E13" "  +2.90%  "
Set-TextValue "E14" "  +4.99%  "
Set-TextValue "D15" "3.176.66"
Set-TextValue "E15" "  +5.40%  "
Set-TextValue "D16" "2.750.16"
Set-TextValue "E16" "  +5.44%  "
Set-TextValue "D17" "0.877"
Set-TextValue "E17" "  +1.93%  "
Set-TextValue "D18" "51.558.53"
Set-TextValue "E18" "  +6.19%  "
Set-TextValue "D19" "3.12"
Set-TextValue "E19" "  +5.92%  "
Set-TextValue "D20" "13.47"
Set-TextValue "E20" "  +5.32%  "
Set-TextValue "E21" "  +2.49%  "
Set-TextValue "D22" "0.0₃0976"
Set-TextValue "E22" "  +3.78%  "
Set-TextValue "D23" "278.60"
Set-TextValue "E23" "  +3.53%  "
Set-TextValue "E24" "  +1.38%  "
Set-TextValue "E25" "  +5.14%  "
Set-TextValue "D26" "26.70"
Set-TextValue "E27" "  -0.01%  "
Set-TextValue "E28" "  +1.84%  "
Set-TextValue "E29" "  +0.34%  "
Set-TextValue "E30" "  +2.03%  "
Set-TextValue "D31" "35.00"
Set-TextValue "D32" "49.97"
Set-TextValue "E32" "  +1.55%  "
Set-TextValue "D33" "5.55"
Set-TextValue "E33" "  +1.96%  "
Set-TextValue "E34" "  +3.00%  "
Set-TextValue "E35" "  -0.03%  "
Set-TextValue "D36" "18.91"
Set-TextValue "E36" "  -0.16%  "
Set-TextValue "D37" "4.96"
Set-TextValue "E37" "  -0.14%  "
Set-TextValue "D38" "2.06"
Set-TextValue "E38" "  +2.15%  "
Set-TextValue "E39" "  +1.38%  "
Set-TextValue "D40" "127.25"
Set-TextValue "E40" "  +0.98%  "
Set-TextValue "E41" "  +9.42%  "
Set-TextValue "D42" "22.94"
Set-TextValue "E42" "  +4.06%  "
Set-TextValue "E43" "  +2.45%  "
Set-TextValue "E44" "  +6.57%  "
Set-TextValue "D45" "2.38"
Set-TextValue "E45" "  +12.56%  "
Set-TextValue "D46" "2.090.61"
Set-TextValue "E46" "  +1.53%  "
Set-TextValue "D47" "3.32"
Set-TextValue "E47" "  +3.15%  "
Set-TextValue "E48" "  +3.58%  "
Set-TextValue "D49" "5.52"
Set-TextValue "E49" "  +7.58%  "
Set-TextValue "D50" "8.92"
Set-TextValue "E50" "  +0.54%  "
Set-TextValue "D51" "59.67"
